$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, matching the workbook's original
# inline-string cells (column D prices / E volume%), so numeric-looking
# strings like "1.00" or "0.507" are not silently coerced into numbers.
# Forcing the Text number format before the assignment keeps the type,
# and resetting the style back to Normal afterwards avoids leaving a
# stray number-format style on the cell.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextValue "D2" "63.094.16"
Set-TextValue "E2" "  +2.89%  "

# --- Row 3: Ethereum ---
Set-TextValue "D3" "2.951.68"
Set-TextValue "E3" "  +0.84%  "

# --- Row 4: TetherUSD ---
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.07%  "

# --- Row 5: BNB ---
Set-TextValue "D5" "594.76"
Set-TextValue "E5" "  -0.43%  "

# --- Row 6: Solana ---
Set-TextValue "D6" "148.28"
Set-TextValue "E6" "  +2.57%  "

# --- Row 7: USDC ---
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.03%  "

# --- Row 8: LidoStakedEther ---
Set-TextValue "D8" "2.949.25"
Set-TextValue "E8" "  +0.80%  "

# --- Row 9: XRP ---
Set-TextValue "D9" "0.507"
Set-TextValue "E9" "  +1.39%  "

# --- Row 10: Toncoin ---
Set-TextValue "D10" "7.12"
Set-TextValue "E10" "  +2.85%  "

# --- Row 11: Dogecoin ---
Set-TextValue "D11" "0.150"
Set-TextValue "E11" "  +6.67%  "

# --- Row 12: Cardano ---
Set-TextValue "D12" "0.440"
Set-TextValue "E12" "  +0.58%  "

# --- Row 13: ShibaInu ---
Set-TextValue "D13" "0.0000235"
Set-TextValue "E13" "  +5.05%  "

# --- Row 14: Avalanche ---
Set-TextValue "D14" "32.76"
Set-TextValue "E14" "  -1.98%  "

# --- Row 15: TRON (price unchanged) ---
Set-TextValue "E15" "  -0.63%  "

# --- Row 16: WrappedliquidstakedEther2.0 ---
Set-TextValue "D16" "3.442.46"
Set-TextValue "E16" "  +0.96%  "

# --- Row 17: WrappedBTC ---
Set-TextValue "D17" "63.060.83"
Set-TextValue "E17" "  +2.86%  "

# --- Row 18: Polkadot ---
Set-TextValue "D18" "6.69"
Set-TextValue "E18" "  +0.19%  "

# --- Row 19: WrappedEther ---
Set-TextValue "D19" "2.957.42"
Set-TextValue "E19" "  +1.09%  "

# --- Row 20: BitcoinCash ---
Set-TextValue "D20" "442.51"
Set-TextValue "E20" "  +2.57%  "

# --- Row 21: Chainlink ---
Set-TextValue "D21" "13.48"
Set-TextValue "E21" "  -0.19%  "

# --- Row 22: Polygon ---
Set-TextValue "D22" "0.667"
Set-TextValue "E22" "  -0.96%  "

# --- Row 23: Uniswap ---
Set-TextValue "D23" "7.01"
Set-TextValue "E23" "  -0.94%  "

# --- Row 24/25: RenderToken & Litecoin swap places (with updated prices) ---
Set-TextValue "B24" "Litecoin"
Set-TextValue "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "80.99"
Set-TextValue "E24" "  -1.02%  "

Set-TextValue "B25" "RenderToken"
Set-TextValue "C25" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D25" "11.14"
Set-TextValue "E25" "  +2.60%  "

# --- Row 26: Fetch.AI ---
Set-TextValue "D26" "2.13"
Set-TextValue "E26" "  -2.20%  "

# --- Row 27: InternetComputer(DFINITY) ---
Set-TextValue "D27" "11.76"
Set-TextValue "E27" "  +0.39%  "

# --- Row 28: Dai (price unchanged) ---
Set-TextValue "E28" "  -0.02%  "

# --- Row 29: NEARProtocol ---
Set-TextValue "D29" "7.27"
Set-TextValue "E29" "  +5.82%  "

# --- Row 30: ImmutableX ---
Set-TextValue "D30" "2.19"
Set-TextValue "E30" "  -0.45%  "

# --- Row 31: PancakeSwap (price unchanged) ---
Set-TextValue "E31" "  +0.29%  "

# --- Row 32: PEPE ---
Set-TextValue "D32" "0.0000102"
Set-TextValue "E32" "  +15.84%  "

# --- Row 33: EthereumClassic ---
Set-TextValue "D33" "26.46"
Set-TextValue "E33" "  -0.59%  "

# --- Row 34: Hedera ---
Set-TextValue "D34" "0.108"
Set-TextValue "E34" "  -1.01%  "

# --- Row 35: FirstDigitalUSD ---
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  -0.05%  "

# --- Row 36: Mantle ---
Set-TextValue "D36" "0.991"
Set-TextValue "E36" "  -1.64%  "

# --- Row 37: dogwifhat ---
Set-TextValue "D37" "3.10"
Set-TextValue "E37" "  +4.36%  "

# --- Row 38: Filecoin ---
Set-TextValue "D38" "5.60"
Set-TextValue "E38" "  -0.43%  "

# --- Row 39/40: Stacks & OKB swap places (with updated prices) ---
Set-TextValue "B39" "OKB"
Set-TextValue "C39" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D39" "49.66"
Set-TextValue "E39" "  -0.47%  "

Set-TextValue "B40" "Stacks"
Set-TextValue "C40" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "2.04"
Set-TextValue "E40" "  +2.56%  "

# --- Row 41: Cosmos ---
Set-TextValue "D41" "8.50"
Set-TextValue "E41" "  -0.50%  "

# --- Row 42: Kaspa ---
Set-TextValue "D42" "0.117"
Set-TextValue "E42" "  -4.07%  "

# --- Row 43: TheGraph ---
Set-TextValue "D43" "0.281"
Set-TextValue "E43" "  +0.25%  "

# --- Row 44: Arweave ---
Set-TextValue "D44" "38.71"
Set-TextValue "E44" "  -8.08%  "

# --- Row 45: Monero ---
Set-TextValue "D45" "135.55"
Set-TextValue "E45" "  +1.46%  "

# --- Row 46: Maker ---
Set-TextValue "D46" "2.692.35"
Set-TextValue "E46" "  -0.28%  "

# --- Row 47: VeChain (price unchanged) ---
Set-TextValue "E47" "  -2.08%  "

# --- Row 48: Bittensor ---
Set-TextValue "D48" "359.78"
Set-TextValue "E48" "  -1.62%  "

# --- Row 50: Stellar (price unchanged) ---
Set-TextValue "E50" "  -0.48%  "

# --- Row 51: InjectiveProtocol ---
Set-TextValue "D51" "22.78"
Set-TextValue "E51" "  -3.08%  "
